$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 488.44116
$ws.Range("I135").Value = 457.03333
$ws.Range("J135").Value = 724
$ws.Range("K135").Value = 4113.29997
$ws.Range("L135").Value = 6516
$ws.Range("M135").Value = -1578.29997
$ws.Range("N135").Value = -11586
$ws.Range("H137").Value = 1248.8445
$ws.Range("I137").Value = 1106.8387
$ws.Range("J137").Value = 1563.2858
$ws.Range("K137").Value = 3320.5161
$ws.Range("L137").Value = 4689.857400000001
$ws.Range("M137").Value = -770.5160999999998
$ws.Range("N137").Value = -9789.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1599.909
$ws.Range("I63").Value = 1879.8
$ws.Range("J63").Value = 1366.6666
$ws.Range("K63").Value = 1879.8
$ws.Range("L63").Value = 1366.6666
$ws.Range("M63").Value = -1193.8
$ws.Range("N63").Value = -2738.6666
$ws.Range("H66").Value = 1599.909
$ws.Range("I66").Value = 1879.8
$ws.Range("J66").Value = 1366.6666
$ws.Range("K66").Value = 9399
$ws.Range("L66").Value = 6833.333000000001
$ws.Range("M66").Value = -5967
$ws.Range("N66").Value = -13697.333
$ws.Range("H122").Value = 1720.4375
$ws.Range("I122").Value = 1422.5518
$ws.Range("K122").Value = 4267.6554
$ws.Range("M122").Value = -1817.6554
$ws.Range("H132").Value = 1845.0588
$ws.Range("I132").Value = 921.875
$ws.Range("J132").Value = 3399.8948
$ws.Range("K132").Value = 2765.625
$ws.Range("L132").Value = 10199.6844
$ws.Range("M132").Value = -235.625
$ws.Range("N132").Value = -15259.6844

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H31").Value = 14330
$ws.Range("I31").Value = 3000
$ws.Range("J31").Value = 19995
$ws.Range("K31").Value = 3000
$ws.Range("L31").Value = 19995
$ws.Range("M31").Value = -2748
$ws.Range("N31").Value = -20499
$ws.Range("H107").Value = 740.05884
$ws.Range("I107").Value = 665.5
$ws.Range("J107").Value = 846.5714
$ws.Range("K107").Value = 665.5
$ws.Range("L107").Value = 846.5714
$ws.Range("M107").Value = 1254.5
$ws.Range("N107").Value = -4686.5714
$ws.Range("H134").Value = 1229.8286
$ws.Range("I134").Value = 906.4
$ws.Range("J134").Value = 2415.7334
$ws.Range("K134").Value = 2719.2
$ws.Range("L134").Value = 7247.2002
$ws.Range("M134").Value = -184.1999999999998
$ws.Range("N134").Value = -12317.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4833.564
$ws.Range("I99").Value = 1880.3182
$ws.Range("K99").Value = 1880.3182
$ws.Range("M99").Value = -382.3181999999999
$ws.Range("H126").Value = 4833.564
$ws.Range("I126").Value = 1880.3182
$ws.Range("K126").Value = 5640.9546
$ws.Range("M126").Value = -3170.9546
$ws.Range("H132").Value = 2415.3
$ws.Range("I132").Value = 1192.3334
$ws.Range("J132").Value = 4249.75
$ws.Range("K132").Value = 3577.0002
$ws.Range("L132").Value = 12749.25
$ws.Range("M132").Value = -1047.0002
$ws.Range("N132").Value = -17809.25
$ws.Range("H134").Value = 1716.02
$ws.Range("I134").Value = 1504.3489
$ws.Range("K134").Value = 4513.0467
$ws.Range("M134").Value = -1978.0467

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1576.0667
$ws.Range("I50").Value = 1313
$ws.Range("J50").Value = 1806.25
$ws.Range("K50").Value = 3939
$ws.Range("L50").Value = 5418.75
$ws.Range("M50").Value = -3458
$ws.Range("N50").Value = -6380.75
$ws.Range("H53").Value = 1576.0667
$ws.Range("I53").Value = 1313
$ws.Range("J53").Value = 1806.25
$ws.Range("K53").Value = 3939
$ws.Range("L53").Value = 5418.75
$ws.Range("M53").Value = -3458
$ws.Range("N53").Value = -6380.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5652330
$ws.Range("I11").Value = 6564225
$ws.Range("J11").Value = 2004750
$ws.Range("K11").Value = 6564225
$ws.Range("L11").Value = 2004750
$ws.Range("M11").Value = -6564086
$ws.Range("N11").Value = -2005028
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10586
$ws.Range("H21").Value = 227263.64
$ws.Range("I21").Value = 100000
$ws.Range("J21").Value = 333316.66
$ws.Range("K21").Value = 100000
$ws.Range("L21").Value = 333316.66
$ws.Range("M21").Value = -99827
$ws.Range("N21").Value = -333662.66
$ws.Range("H30").Value = 227263.64
$ws.Range("I30").Value = 100000
$ws.Range("J30").Value = 333316.66
$ws.Range("K30").Value = 100000
$ws.Range("L30").Value = 333316.66
$ws.Range("M30").Value = -99895
$ws.Range("N30").Value = -333526.66
$ws.Range("H126").Value = 1825.6207
$ws.Range("I126").Value = 1746.7727
$ws.Range("J126").Value = 2073.4285
$ws.Range("K126").Value = 5240.3181
$ws.Range("L126").Value = 6220.2855
$ws.Range("M126").Value = -2770.3181
$ws.Range("N126").Value = -11160.2855
$ws.Range("H132").Value = 2261.1562
$ws.Range("I132").Value = 2230.5557
$ws.Range("J132").Value = 2300.5
$ws.Range("K132").Value = 6691.6671
$ws.Range("L132").Value = 6901.5
$ws.Range("M132").Value = -4161.6671
$ws.Range("N132").Value = -11961.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 7966.6665
$ws.Range("I23").Value = 7966.6665
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 7966.6665
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -7736.6665
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 50000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 50000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 50000
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -50590
$ws.Range("H30").Value = 700
$ws.Range("I30").Value = 700
$ws.Range("K30").Value = 700
$ws.Range("M30").Value = -592
$ws.Range("H122").Value = 7939763.5
$ws.Range("I122").Value = 10102637
$ws.Range("J122").Value = 9226.666999999999
$ws.Range("K122").Value = 30307911
$ws.Range("L122").Value = 27680.001
$ws.Range("M122").Value = -30305461
$ws.Range("N122").Value = -32580.001
$ws.Range("H136").Value = 8548570
$ws.Range("I136").Value = 1618.1892
$ws.Range("J136").Value = 166667170
$ws.Range("K136").Value = 4854.5676
$ws.Range("L136").Value = 500001510
$ws.Range("M136").Value = -2304.5676
$ws.Range("N136").Value = -500006610

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 2510830
$ws.Range("H132").Value = 1724.7966
$ws.Range("I132").Value = 1369.902
$ws.Range("J132").Value = 3987.25
$ws.Range("K132").Value = 4109.706
$ws.Range("L132").Value = 11961.75
$ws.Range("M132").Value = -1579.706
$ws.Range("N132").Value = -17021.75
$ws.Range("H136").Value = 4976596.5
$ws.Range("I136").Value = 7576341
$ws.Range("J136").Value = 3171.4783
$ws.Range("K136").Value = 22729023
$ws.Range("L136").Value = 9514.4349
$ws.Range("M136").Value = -22726473
$ws.Range("N136").Value = -14614.4349
